# Updated cryptos list on Sun Oct 29 10:57:18 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
}

# Row 2 - Bitcoin
Set-Text "D2" "34.213.55"
Set-Text "E2" "  +0.34%  "

# Row 3 - Ethereum
Set-Text "D3" "1.786.31"
Set-Text "E3" "  -0.23%  "

# Row 4 - TetherUSD
Set-Text "E4" "  +0.11%  "

# Row 5 - BNB
Set-Text "D5" "225.90"
Set-Text "E5" "  -0.56%  "

# Row 6 - XRP
Set-Text "E6" "  +1.25%  "

# Row 7 - USDC
Set-Text "E7" "  +0.10%  "

# Row 8 - Solana
Set-Text "E8" "  -0.05%  "

# Row 9 - Cardano
Set-Text "E9" "  -0.02%  "

# Row 10 - Dogecoin
Set-Text "E10" "  -0.63%  "

# Row 11 - TRON
Set-Text "E11" "  +0.76%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-Text "D12" "2.044.80"
Set-Text "E12" "  -0.16%  "

# Row 13 & 14 - WrappedEther and Chainlink swap ranking order
Set-Text "B13" "Chainlink"
Set-Text "C13" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-Text "D13" "10.99"
Set-Text "E13" "  -4.36%  "

Set-Text "B14" "WrappedEther"
Set-Text "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-Text "D14" "1.781.01"
Set-Text "E14" "  -0.61%  "

# Row 15 - Polygon
Set-Text "E15" "  +0.51%  "

# Row 16 - WrappedBTC
Set-Text "D16" "34.195.92"
Set-Text "E16" "  +0.34%  "

# Row 17 - Polkadot
Set-Text "E17" "  -0.25%  "

# Row 18 - Litecoin
Set-Text "D18" "67.96"
Set-Text "E18" "  +0.10%  "

# Row 19 - ShibaInu
Set-Text "D19" "0.0₃0801"
Set-Text "E19" "  +2.56%  "

# Row 20 - BitcoinCash
Set-Text "D20" "245.75"
Set-Text "E20" "  +0.04%  "

# Row 22 - Avalanche
Set-Text "E22" "  -0.08%  "

# Row 24 - Toncoin
Set-Text "E24" "  +0.32%  "

# Row 25 - Monero
Set-Text "D25" "161.71"
Set-Text "E25" "  -0.08%  "

# Row 26 - Cosmos
Set-Text "E26" "  -0.14%  "

# Row 27 - EthereumClassic
Set-Text "D27" "16.32"
Set-Text "E27" "  +0.09%  "

# Row 28 - Stellar
Set-Text "E28" "  +1.33%  "

# Row 29 - BinanceUSD
Set-Text "E29" "  +0.16%  "

# Row 30 - PancakeSwap
Set-Text "E30" "  -0.76%  "

# Row 31 - Hedera
Set-Text "E31" "  +0.03%  "

# Row 32 - Filecoin
Set-Text "E32" "  +2.07%  "

# Row 33 - InternetComputer(DFINITY)
Set-Text "E33" "  +3.85%  "

# Row 34 - LidoDAOToken
Set-Text "E34" "  -1.62%  "

# Row 35 - Maker
Set-Text "D35" "1.437.99"
Set-Text "E35" "  -0.49%  "

# Row 36 - RenderToken
Set-Text "D36" "2.62"
Set-Text "E36" "  +10.96%  "

# Row 37 - ImmutableX
Set-Text "E37" "  +2.40%  "

# Row 38 - TrustWalletToken
Set-Text "E38" "  +1.73%  "

# Row 39 - VeChain
Set-Text "E39" "  -1.15%  "

# Row 40 - Aave
Set-Text "D40" "81.59"
Set-Text "E40" "  +1.34%  "

# Row 41 - HuobiToken
Set-Text "E41" "  +1.48%  "

# Row 42 - InjectiveProtocol
Set-Text "D42" "14.10"
Set-Text "E42" "  +5.80%  "

# Row 43 - MXToken
Set-Text "E43" "  +1.49%  "

# Row 44 - ARBITRUM
Set-Text "E44" "  -0.46%  "

# Row 45 - Kaspa
Set-Text "D45" "0.0520"
Set-Text "E45" "  +2.23%  "

# Row 46 & 47 - FraxShare and WEMIXToken swap ranking order
Set-Text "B46" "WEMIXToken"
Set-Text "C46" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-Text "D46" "1.08"
Set-Text "E46" "  +1.23%  "

Set-Text "B47" "FraxShare"
Set-Text "C47" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-Text "D47" "6.06"
Set-Text "E47" "  -0.20%  "

# Row 48 - RocketPoolETH
Set-Text "D48" "1.940.79"

# Row 49 - Quant
Set-Text "D49" "105.46"
Set-Text "E49" "  -2.08%  "

# Row 50 - PaxDollar
Set-Text "E50" "  +0.13%  "

# Row 51 - BabyDogeCoin
Set-Text "E51" "  -6.43%  "
